# Apply the authored content edits to schedule.xlsx (Sheet1).
# The bulk of the raw XML diff (namespace bumps, fileVersion/rupBuild,
# window pixel sizes, wholesale styles.xml font/fill/border table
# rebuild, theme extraClrSchemeLst, per-row x14ac:dyDescent, etc.) is
# produced automatically by a newer Excel build re-saving the workbook
# and carries no semantic content; it isn't reachable via the Excel
# object model and is intentionally left alone. Below are the genuine
# content edits: four date/time corrections, one cell's text extended
# with newly-appended IDs (as rich text, matching the run split/ font
# seen in the shared-string table), the active selection, and the two
# duplicate-value conditional-format rules losing their absolute ($)
# row anchors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: planning date/time corrections -----------------------
# Rows 2 & 18 and rows 36 / 43 all move forward exactly 24 days (same
# time-of-day); row 12 moves forward 30 minutes.
$ws.Range("A2").Value2  = 45273.895833333336
$ws.Range("A18").Value2 = 45273.895833333336
$ws.Range("A36").Value2 = 45273.5625
$ws.Range("A43").Value2 = 45273.8125
$ws.Range("A12").Value2 = 45253.916666666664

# --- B12: append newly-posted IDs to the existing value --------------
# Original text (shared string) was ",103482326003878,,17841456036806884".
# Three more IDs were appended; the appended tail is written back in a
# distinct (Microsoft YaHei) run, same as the rest of the sheet's body font.
$origPart = ",103482326003878,,17841456036806884"
$newPart  = ",17841461742288388,127520840434805"
$ws.Range("B12").Value2 = $origPart + $newPart
$ws.Range("B12").Font.Name = "Microsoft YaHei"
$ws.Range("B12").Font.Size = 11

$chars = $ws.Range("B12").Characters($origPart.Length + 1, $newPart.Length)
$chars.Font.Name = "Microsoft YaHei"
$chars.Font.Size = 11

# --- Selection / active view -----------------------------------------
$ws.Activate()
$ws.Range("B22").Select()

# --- Conditional formatting: drop the "$" row anchors on both ranges -
# (duplicate-value highlighting on columns C and D). The D-column rule
# is simply re-pointed at the un-anchored range; the C-column rule is
# rebuilt with an identical format so it gets its own dxf record, and
# its priority is put back at 3 (unchanged) so the D rule keeps its
# original priority of 2.
$dRule = $ws.Range("D$1:D$1048576").FormatConditions.Item(1)
$dRule.ModifyAppliesToRange($ws.Range("D1:D1048576"))

$ws.Range("C1:C1048576").FormatConditions.Delete()
$cRule = $ws.Range("C1:C1048576").FormatConditions.AddUniqueValues()
$cRule.DupeUnique = 1
$cRule.Font.Color = 393372
$cRule.Interior.Color = 13551615
$cRule.Priority = 3
